$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: update the date in A1 (+1 day: 2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Step 2: update the unit prices in column D for rows 22-25
$ws.Range("D22").Value = 12264
$ws.Range("D23").Value = 13894
$ws.Range("D24").Value = 18098
$ws.Range("D25").Value = 20048
